# Apply template.docx style-sheet changes:
#  1. Normal style gets spacing-after of 120 twips (6pt).
#  2. Two new custom paragraph styles are added: tei_collation / tei_extent,
#     both based on Heading4, following with Normal, and marked as QuickStyle.

$d = $word.ActiveDocument

# 1. Normal style: add <w:pPr><w:spacing w:after="120"/></w:pPr>
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceAfter = 6

# 2. New style: tei_collation (styleId "teicollation")
$collation = $d.Styles.Add("teicollation", 1)
$collation.NameLocal = "tei_collation"
$collation.BaseStyle = "Heading4"
$collation.NextParagraphStyle = "Normal"
$collation.QuickStyle = $true

# 3. New style: tei_extent (styleId "teiextent")
$extent = $d.Styles.Add("teiextent", 1)
$extent.NameLocal = "tei_extent"
$extent.BaseStyle = "Heading4"
$extent.NextParagraphStyle = "Normal"
$extent.QuickStyle = $true
